# Update the "想去人数" (F column) values for the "展览" and "全部类型" sheets.
# Both sheets contain identical data, and the same set of rows changed in each.

$wb = $excel.ActiveWorkbook

# Row -> new value mapping for column F
$updates = @{
    2  = 255
    3  = 1323
    6  = 223
    8  = 13
    9  = 176
    10 = 127
    11 = 4447
    12 = 6710
    16 = 564
    18 = 4099
    19 = 465
    21 = 47
    22 = 2678
    25 = 164
    26 = 346
    27 = 348
    28 = 394
    29 = 216
    32 = 1014
    33 = 59
    35 = 77
    37 = 495
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
